# Append: 2025-12-07 06:25 JST
# Applies the new scrape pass over the "ランサーズ" (Lancers) listing sheet:
#   - refresh the timestamp column for the rows that stay
#   - rows 2-6 get new listing data (titles/prices/urls/scores/skill tags)
#   - the old rows 7-10 are gone entirely (dimension shrinks to A1:H6)
#   - column H narrows back down from 16 to 12 characters

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2025-12-07 06:25:35"

# --- New data for rows 2-6 --------------------------------------------------
# row, title, price, url, score, skills  (category/deadline columns are unchanged)
$rows = @(
    @{ Row = 2; Title = "【急募】あなたAIクローン構築パートナー募集!"; Price = "500,000 円 ~ 1,000,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5448719"; Score = 310; Skills = "🔥AI,Ai" },
    @{ Row = 3; Title = "【自動化】食べログの「いいね」「フォロー」作業を効率化したい"; Price = "10,000 円 ~ 20,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5448871"; Score = 145; Skills = "◆効率化,自動化" },
    @{ Row = 4; Title = "【システム開発】顧客予約サインシステムの構築依頼"; Price = "20,000 円 ~ 50,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5448905"; Score = 113; Skills = "◆開発,システム開発" },
    @{ Row = 5; Title = "laravelで作った顧客管理Webの表示情報を、他の顧客管理Web上にコピペ入力する作業を自動化"; Price = "20,000 円 ~ 50,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5448875"; Score = 108; Skills = "◆自動化 ◇管理" },
    @{ Row = 6; Title = "【急募】HPリニューアルで業務フロー自動化を実現するプロ募集!"; Price = "100,000 円 ~ 200,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5448817"; Score = 88; Skills = "◆自動化" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $timestamp
    $ws.Cells.Item($row, 2).Value = $r.Title
    $ws.Cells.Item($row, 4).Value = $r.Price
    $ws.Cells.Item($row, 6).Value = $r.Url
    $ws.Cells.Item($row, 7).Value = $r.Score
    $ws.Cells.Item($row, 8).Value = $r.Skills
}

# --- Drop the old rows 7-10 (they no longer appear in this scrape) ---------
$ws.Rows("7:10").Delete()

# --- Hyperlinks: rebuild so only F2:F6 carry links, pointing at the new URLs
# (the engine's Hyperlinks collection only supports clearing everything /
# re-adding, so remove then re-add in order to keep relationship ids stable)
$ws.Hyperlinks.Delete()
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r.Row, 6)
    $ws.Hyperlinks.Add($cell, $r.Url) | Out-Null
    # Hyperlinks.Add() stamps a fresh (duplicate) font style on the cell;
    # pin it back to the workbook's existing "Hyperlink" cell style so F2:F6
    # keep the same style index they already had.
    $cell.Style = "Hyperlink"
}

# --- Column H width back to 12 (character width, matches stored width=12) --
$ws.Columns("H").ColumnWidth = 11.166666666666666

Write-Output "edit complete"
